# Regenerate the handoff report: refresh the "Latest Handoff" timestamp for
# every file that is still outstanding (i.e. every row whose status is not
# "Handed back: in sync with en-US" and not "In Translation") on each sheet.
#
#  - Overview sheet:  column D ("Latest Handoff Date")
#  - zh-cn sheet:      column E ("Latest Handoff Datetime")
#  - de-de sheet:      column E ("Latest Handoff Datetime")

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D4").Value = "2016-59-14 09:59:45"
$overview.Range("D6").Value = "2016-59-14 09:59:45"
$overview.Range("D7").Value = "2016-59-14 09:59:45"
$overview.Range("D8").Value = "2016-59-14 09:59:45"
$overview.Range("D9").Value = "2016-59-14 09:59:45"
$overview.Range("D10").Value = "2016-59-14 09:59:45"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E4").Value = "2016-03-14 09:59:37"
$zhcn.Range("E6").Value = "2016-03-14 09:59:37"
$zhcn.Range("E7").Value = "2016-03-14 09:59:37"
$zhcn.Range("E8").Value = "2016-03-14 09:59:37"
$zhcn.Range("E9").Value = "2016-03-14 09:59:37"
$zhcn.Range("E10").Value = "2016-03-14 09:59:37"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E4").Value = "2016-03-14 09:59:45"
$dede.Range("E6").Value = "2016-03-14 09:59:45"
$dede.Range("E7").Value = "2016-03-14 09:59:45"
$dede.Range("E8").Value = "2016-03-14 09:59:45"
$dede.Range("E9").Value = "2016-03-14 09:59:45"
$dede.Range("E10").Value = "2016-03-14 09:59:45"
